{"js": "// Helper: replace the contents of a Word.Range with explicit OOXML\n// (paragraph(s) of runs / proofErr markers / bookmarks), using the\n// pkg:package wrapper required by Range.insertOoxml.\nfunction wrapOoxml(bodyXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    bodyXml +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the paragraphs we need purely by their current text so the\n// script is resilient to exact indices.\nconst texts = paragraphs.items.map((p) => p.text);\n\nfunction findIndex(expected) {\n  const idx = texts.indexOf(expected);\n  if (idx === -1) {\n    throw new Error(\"Could not find paragraph with text: \" + expected);\n  }\n  return idx;\n}\n\nconst idxHomepage = findIndex(\"Homepage: latest news , \");\nconst idxNavBar = findIndex(\"Fix navigationbar add background image \");\nconst idxWhatLeft = findIndex(\"What left: \");\nconst idxEducation = findIndex(\"2- education and biography read and display \");\nconst idxResearch = findIndex(\"3-research and publication read and display \");\n\n// 1) \"Homepage: latest news , \" -> split into 3 runs with a gramStart/gramEnd\n//    proofErr pair wrapping \"news ,\".\nparagraphs.items[idxHomepage].getRange().insertOoxml(\n  wrapOoxml(\n    \"<w:p>\" +\n      '<w:r><w:t xml:space=\"preserve\">Homepage: latest </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      \"<w:r><w:t>news ,</w:t></w:r>\" +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n      \"</w:p>\"\n  ),\n  Word.InsertLocation.replace\n);\n\n// 2) \"Fix navigationbar add background image \" -> split into 3 runs with a\n//    spellStart/spellEnd proofErr pair wrapping \"navigationbar\".\nparagraphs.items[idxNavBar].getRange().insertOoxml(\n  wrapOoxml(\n    \"<w:p>\" +\n      '<w:r><w:t xml:space=\"preserve\">Fix </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      \"<w:r><w:t>navigationbar</w:t></w:r>\" +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> add background image </w:t></w:r>' +\n      \"</w:p>\"\n  ),\n  Word.InsertLocation.replace\n);\n\n// 3) \"What left: \" -> split into 3 runs with a gramStart/gramEnd proofErr\n//    pair wrapping \"left:\".\nparagraphs.items[idxWhatLeft].getRange().insertOoxml(\n  wrapOoxml(\n    \"<w:p>\" +\n      '<w:r><w:t xml:space=\"preserve\">What </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      \"<w:r><w:t>left:</w:t></w:r>\" +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n      \"</w:p>\"\n  ),\n  Word.InsertLocation.replace\n);\n\n// 4) The \"2- education ...\" and \"3-research ...\" paragraphs are replaced\n//    together: the first gets a gramStart/gramEnd proofErr pair around\n//    \"education\"; the second drops its trailing space and its bookmark\n//    (the _GoBack bookmark is relocated further below). New paragraphs\n//    (including three blank ones and four new content paragraphs) are\n//    appended right after, in one combined OOXML insert so the bookmark\n//    move and paragraph insert happen atomically.\nconst eduRange = paragraphs.items[idxEducation].getRange();\nconst researchRange = paragraphs.items[idxResearch].getRange();\nconst combinedRange = eduRange.expandTo(researchRange);\n\ncombinedRange.insertOoxml(\n  wrapOoxml(\n    \"<w:p>\" +\n      '<w:r><w:t xml:space=\"preserve\">2- </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      \"<w:r><w:t>education</w:t></w:r>\" +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> and biography read and display </w:t></w:r>' +\n      \"</w:p>\" +\n      \"<w:p>\" +\n      \"<w:r><w:t>3-research and publication read and display</w:t></w:r>\" +\n      \"</w:p>\" +\n      \"<w:p/>\" +\n      \"<w:p/>\" +\n      \"<w:p/>\" +\n      \"<w:p>\" +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      \"<w:r><w:t>Gallary</w:t></w:r>\" +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\">: header= pictures and videos </w:t></w:r>' +\n      \"</w:p>\" +\n      \"<w:p>\" +\n      '<w:r><w:t xml:space=\"preserve\">Navigation bar: make it smaller and change the </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      \"<w:r><w:t>colo</w:t></w:r>\" +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      \"</w:p>\" +\n      \"<w:p>\" +\n      \"<w:r><w:t>Home: slide show being able to change im</w:t></w:r>\" +\n      '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n      '<w:bookmarkEnd w:id=\"0\"/>' +\n      '<w:r><w:t xml:space=\"preserve\">ages </w:t></w:r>' +\n      '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n      \"</w:p>\" +\n      \"<w:p/>\" +\n      \"<w:p>\" +\n      '<w:r><w:t xml:space=\"preserve\">Upload website on online server </w:t></w:r>' +\n      \"</w:p>\"\n  ),\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Apply proofing-mark run splits (gramStart/gramEnd, spellStart/spellEnd)\n# to four existing paragraphs, tighten the \"3-research...\" paragraph\n# (drop trailing space, drop the _GoBack bookmark there), and append a\n# batch of new paragraphs (three blank + four with content, one of which\n# now hosts the relocated _GoBack bookmark) right after it.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($doc, $text) {\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $doc.Paragraphs($i)\n        if ($p.Range.Text.TrimEnd(\"`r\") -eq $text) {\n            return $i\n        }\n    }\n    throw \"Paragraph not found: $text\"\n}\n\nfunction New-PkgXml($bodyXml) {\n    return '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $bodyXml + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\n# 1) \"Homepage: latest news , \" -> split into 3 runs with a gramStart/gramEnd\n#    proofErr pair wrapping \"news ,\".\n$idx = Find-ParagraphIndex $d \"Homepage: latest news , \"\n$p = $d.Paragraphs($idx)\n$xml = New-PkgXml(\n    '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">Homepage: latest </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>news ,</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '</w:p>'\n)\n$p.Range.InsertXML($xml)\n\n# 2) \"Fix navigationbar add background image \" -> split into 3 runs with a\n#    spellStart/spellEnd proofErr pair wrapping \"navigationbar\".\n$idx = Find-ParagraphIndex $d \"Fix navigationbar add background image \"\n$p = $d.Paragraphs($idx)\n$xml = New-PkgXml(\n    '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">Fix </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>navigationbar</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> add background image </w:t></w:r>' +\n    '</w:p>'\n)\n$p.Range.InsertXML($xml)\n\n# 3) \"What left: \" -> split into 3 runs with a gramStart/gramEnd proofErr\n#    pair wrapping \"left:\".\n$idx = Find-ParagraphIndex $d \"What left: \"\n$p = $d.Paragraphs($idx)\n$xml = New-PkgXml(\n    '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">What </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>left:</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '</w:p>'\n)\n$p.Range.InsertXML($xml)\n\n# 4) The \"2- education ...\" and \"3-research ...\" paragraphs are replaced\n#    together: the first gets a gramStart/gramEnd proofErr pair around\n#    \"education\"; the second drops its trailing space and its bookmark\n#    (the _GoBack bookmark is relocated further below). New paragraphs\n#    (three blank ones plus four new content paragraphs) are appended\n#    right after, in one combined InsertXML call so the bookmark move and\n#    the paragraph insert happen atomically.\n$idxEdu = Find-ParagraphIndex $d \"2- education and biography read and display \"\n$idxRes = Find-ParagraphIndex $d \"3-research and publication read and display \"\n$p1 = $d.Paragraphs($idxEdu)\n$p2 = $d.Paragraphs($idxRes)\n$combined = $d.Range($p1.Range.Start, $p2.Range.End)\n\n$body = ''\n$body += '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">2- </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>education</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> and biography read and display </w:t></w:r>' +\n    '</w:p>'\n$body += '<w:p><w:r><w:t>3-research and publication read and display</w:t></w:r></w:p>'\n$body += '<w:p/><w:p/><w:p/>'\n$body += '<w:p>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Gallary</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">: header= pictures and videos </w:t></w:r>' +\n    '</w:p>'\n$body += '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">Navigation bar: make it smaller and change the </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>colo</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '</w:p>'\n$body += '<w:p>' +\n    '<w:r><w:t>Home: slide show being able to change im</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">ages </w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '</w:p>'\n$body += '<w:p/>'\n$body += '<w:p><w:r><w:t xml:space=\"preserve\">Upload website on online server </w:t></w:r></w:p>'\n\n$xml = New-PkgXml($body)\n$combined.InsertXML($xml)\n"}
